$d = $word.ActiveDocument

# --- Change 1: expand the "regression test program" bullet about starting
# data/expected outputs into six detailed validation bullets. ---
$rng = $d.Content
$rng.Find.Execute("Identify a set of starting data and a set of expected outputs from running the program", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$para = $rng.Paragraphs(1)

$xmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$newBulletsXml = @"
<w:p $xmlNs><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Validate that the collection of files in the raw folder is in a date folder and contains fits, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>png</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> for thumbnail, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>png</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> for full image and (one more file?)</w:t></w:r></w:p>
<w:p $xmlNs><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Validate that the collection of files in the raw folder has the correct entries in the fits table and the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>fits_by_target</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> table.</w:t></w:r></w:p>
<w:p $xmlNs><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Validate that the bias files from a collection for a date in the raw folder are copied to the bias folder.</w:t></w:r></w:p>
<w:p $xmlNs><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Validate that a valid master bias fits file is created from a combination of </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>all of</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> the bias frames for a date.</w:t></w:r></w:p>
<w:p $xmlNs><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Validate that the fits and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>fits_by_target</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> table entries are present for the bias frames and the master bias frame.</w:t></w:r></w:p>
<w:p $xmlNs><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>&#8230;</w:t></w:r></w:p>
"@

$para.Range.InsertXML($newBulletsXml)

# --- Change 2: drop the stray lastRenderedPageBreak that used to sit before
# "Establish a code push and test process" (it moved to the bullet above). ---
$rng2 = $d.Content
$rng2.Find.Execute("Establish a code push and test process", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pushPara = $rng2.Paragraphs(1)
$pushXml = @"
<w:p $xmlNs><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Establish a code push and test process</w:t></w:r></w:p>
"@
$pushPara.Range.InsertXML($pushXml)

# --- Change 3: split the "reset" sentence so it carries a grammar-check
# marker around the word "reset". ---
$rng3 = $d.Content
$rng3.Find.Execute("Once fixes are in place, the test environment should be reset and the regression test run again.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$resetPara = $rng3.Paragraphs(1)
$resetXml = @"
<w:p $xmlNs><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Once fixes are in place, the test environment should be </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>reset</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> and the regression test run again.</w:t></w:r></w:p>
"@
$resetPara.Range.InsertXML($resetXml)
